$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New bug report added as row 51 (row 50 intentionally left blank, consistent with
# other gaps already present in this bug log).
$ws.Range("A51").Value = "SB"

$ws.Range("B51").Value = 44125
$ws.Range("D51").Value = 44125
# Reuse the existing date-formatted style (the same one already used by B49/D49)
# rather than letting a brand-new number format style get created.
$ws.Range("B49").Copy()
$ws.Range("B51").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D51").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("C51").Value = 'in sumby if the factor labels have spaces in then the heading in the output is "tidied". Disable this. '
$ws.Range("C51").Style = "Good"

# Move the view down to the newly added row, mirroring the scrolled/selected state
# the workbook was left in after the edit.
$ws.Application.ActiveWindow.ScrollRow = 40
$ws.Range("A51").Select()
